# Auto-generated edit script: updates crypto price/volume table cells
# to match the refreshed data pulled on Thu Jun 13 15:34:53 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "66.981.43"
$ws.Range("E2").Formula = "  -4.08%  "
$ws.Range("D3").Formula = "3.478.53"
$ws.Range("E3").Formula = "  -4.12%  "
$ws.Range("E4").Formula = "  +0.15%  "
$ws.Range("D5").Formula = "'601.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "  -4.66%  "
$ws.Range("D6").Formula = "'147.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "  -7.30%  "
$ws.Range("D7").Formula = "3.478.70"
$ws.Range("E7").Formula = "  -4.00%  "
$ws.Range("E8").Formula = "  +0.01%  "
$ws.Range("D9").Formula = "'0.484"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "  -2.67%  "
$ws.Range("E10").Formula = "  -4.32%  "
$ws.Range("D11").Formula = "'7.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "  +1.45%  "
$ws.Range("D12").Formula = "'0.425"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "  -3.71%  "
$ws.Range("D13").Formula = "'0.0000213"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "  -7.11%  "
$ws.Range("D14").Formula = "4.066.93"
$ws.Range("E14").Formula = "  -4.14%  "
$ws.Range("D15").Formula = "'31.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "  -5.51%  "
$ws.Range("D16").Formula = "3.480.91"
$ws.Range("E16").Formula = "  -3.82%  "
$ws.Range("D17").Formula = "67.089.96"
$ws.Range("E17").Formula = "  -4.08%  "
$ws.Range("E18").Formula = "  -0.73%  "
$ws.Range("D19").Formula = "'6.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "  -3.90%  "
$ws.Range("D20").Formula = "'15.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "  -4.60%  "
$ws.Range("D21").Formula = "'9.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "  -2.83%  "
$ws.Range("D22").Formula = "'441.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "  -5.00%  "
$ws.Range("D23").Formula = "'0.619"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "  -4.41%  "
$ws.Range("D24").Formula = "'78.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "  -0.23%  "
$ws.Range("B25").Formula = "Dai"
$ws.Range("C25").Formula = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Formula = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "  -0.03%  "
$ws.Range("B26").Formula = "WrappedeETH"
$ws.Range("C26").Formula = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Formula = "3.623.55"
$ws.Range("E26").Formula = "  -3.96%  "
$ws.Range("E27").Formula = "  -5.54%  "
$ws.Range("D28").Formula = "'0.0000119"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Formula = "  -12.13%  "
$ws.Range("D29").Formula = "'9.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Formula = "  -8.55%  "
$ws.Range("D30").Formula = "'8.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "  -9.55%  "
$ws.Range("D31").Formula = "'2.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "  -6.08%  "
$ws.Range("D32").Formula = "'1.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "  -6.83%  "
$ws.Range("E33").Formula = "  +0.17%  "
$ws.Range("D34").Formula = "'0.165"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Formula = "  -8.23%  "
$ws.Range("D35").Formula = "'25.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Formula = "  -4.57%  "
$ws.Range("D36").Formula = "'6.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "  -7.54%  "
$ws.Range("D37").Formula = "3.474.18"
$ws.Range("E37").Formula = "  -4.21%  "
$ws.Range("E38").Formula = "  -8.52%  "
$ws.Range("B39").Formula = "USDe"
$ws.Range("C39").Formula = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").Formula = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Formula = "  -0.01%  "
$ws.Range("B40").Formula = "Aptos"
$ws.Range("C40").Formula = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Formula = "'7.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "  -7.13%  "
$ws.Range("E42").Formula = "  -8.90%  "
$ws.Range("D43").Formula = "'173.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "  -3.33%  "
$ws.Range("D44").Formula = "'0.0889"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "  -4.53%  "
$ws.Range("D45").Formula = "'5.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "  -5.84%  "
$ws.Range("D46").Formula = "'0.892"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Formula = "  -2.50%  "
$ws.Range("D47").Formula = "'29.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "  -8.63%  "
$ws.Range("D48").Formula = "'46.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "  +0.33%  "
$ws.Range("D49").Formula = "'1.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "  -8.98%  "
$ws.Range("D50").Formula = "'7.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "  -4.18%  "
$ws.Range("D51").Formula = "'2.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "  -12.15%  "
